# Update cryptos worksheet cell values to match the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.224.76"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.075.13"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.24"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.674"
$ws.Range("E6").Value = "  +3.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.63"
$ws.Range("E7").Value = "  +25.24%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.98"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("E10").Value = "  +5.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +8.64%  "
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.75"
$ws.Range("E13").Value = "  +4.04%  "
$ws.Range("D14").Value = "2.345.14"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.826"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  +7.32%  "
$ws.Range("D17").Value = "2.083.87"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "37.170.39"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.84"
$ws.Range("E19").Value = "  +4.15%  "
$ws.Range("D20").Value = "0.0₃0927"
$ws.Range("E20").Value = "  +13.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.14"
$ws.Range("E21").Value = "  +14.91%  "
$ws.Range("E22").Value = "  +6.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.03"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.51"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.39"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.74"
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0637"
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("E34").Value = "  +10.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0891"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -3.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("E39").Value = "  +23.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.85"
$ws.Range("E41").Value = "  +7.78%  "
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.56"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.39"
$ws.Range("E45").Value = "  +26.17%  "
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.57"
$ws.Range("E47").Value = "  +14.51%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.54"
$ws.Range("E48").Value = "  +15.57%  "
$ws.Range("D49").Value = "1.308.53"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.95"
$ws.Range("E51").Value = "  +1.87%  "
